# Adding 4 search test cases to the "Test Cases" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Stable "donor" cells used purely to copy an existing cell style (format-only
# paste) onto new/changed cells. None of these donor cells are themselves
# touched anywhere else in this script, so they stay put as a style source
# throughout the whole run.
#   style 1  -> A42
#   style 4  -> C64
#   style 8  -> C59
#   style 18 -> D2

# --- Fix up existing rows 60, 64, 65 (status flips to SKIP, style alignment on D) ---

# Row 60: result PASS -> SKIP
$ws.Range("E60").Value = "SKIP"

# Row 64: D64 style needs to line up with the rest of column D (style 18); E64 PASS -> SKIP
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D64").PasteSpecial(-4122) | Out-Null
$ws.Range("E64").Value = "SKIP"

# Row 65: D65 style needs to line up with the rest of column D (style 18); E65 PASS -> SKIP
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D65").PasteSpecial(-4122) | Out-Null
$ws.Range("E65").Value = "SKIP"

# --- New row 66 : VerifyPublishPostDisplayed ---
$ws.Range("A42").Copy() | Out-Null
$ws.Range("A66").PasteSpecial(-4122) | Out-Null
$ws.Range("A42").Copy() | Out-Null
$ws.Range("B66").PasteSpecial(-4122) | Out-Null
$ws.Range("C64").Copy() | Out-Null
$ws.Range("C66").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D66").PasteSpecial(-4122) | Out-Null
$ws.Range("A42").Copy() | Out-Null
$ws.Range("E66").PasteSpecial(-4122) | Out-Null

$ws.Range("A66").Value = "VerifyPublishPostDisplayed"
$ws.Range("B66").Value = "OPQA-1190"
$ws.Range("C66").Value = "Verify that Publish a Post option is displayed in Home page and all Record view `npages like Article,Post ,Patent"
$ws.Range("D66").Value = "Y"
$ws.Range("E66").Value = "SKIP"
$ws.Rows.Item(66).RowHeight = 45

# --- New row 67 : VerifyDraftPostTabDisplayForZeroDrafts ---
$ws.Range("A42").Copy() | Out-Null
$ws.Range("A67").PasteSpecial(-4122) | Out-Null
$ws.Range("C59").Copy() | Out-Null
$ws.Range("B67").PasteSpecial(-4122) | Out-Null
$ws.Range("A42").Copy() | Out-Null
$ws.Range("C67").PasteSpecial(-4122) | Out-Null
$ws.Range("C59").Copy() | Out-Null
$ws.Range("D67").PasteSpecial(-4122) | Out-Null
$ws.Range("A42").Copy() | Out-Null
$ws.Range("E67").PasteSpecial(-4122) | Out-Null

# Values written in authored order: description, then jira id, then test name.
$ws.Range("C67").Value = "Verfiy that the Drafts Post tab is not displayed when there are no draft posts"
$ws.Range("B67").Value = "OPQA-1198"
$ws.Range("A67").Value = "VerifyDraftPostTabDisplayForZeroDrafts"
$ws.Range("D67").Value = "Y"
$ws.Range("E67").Value = "PASS"

# --- New row 68 : DeleteDraftPostFromPostModal ---
$ws.Range("A42").Copy() | Out-Null
$ws.Range("A68").PasteSpecial(-4122) | Out-Null
$ws.Range("A42").Copy() | Out-Null
$ws.Range("B68").PasteSpecial(-4122) | Out-Null
$ws.Range("A42").Copy() | Out-Null
$ws.Range("C68").PasteSpecial(-4122) | Out-Null
$ws.Range("A42").Copy() | Out-Null
$ws.Range("D68").PasteSpecial(-4122) | Out-Null
$ws.Range("A42").Copy() | Out-Null
$ws.Range("E68").PasteSpecial(-4122) | Out-Null

$ws.Range("A68").Value = "DeleteDraftPostFromPostModal"
$ws.Range("B68").Value = "OPQA-1200"
$ws.Range("C68").Value = "Verfiy that user is able to delete the draft post from post modal"
$ws.Range("D68").Value = "Y"
$ws.Range("E68").Value = "PASS"

# --- View state: selection moves to A73 ---
$ws.Range("A73").Select() | Out-Null
